# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.932.65"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "'1.894.97"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'0.7758"
$ws.Range("E5").Value = "  -2.90%  "

$ws.Range("D6").Value = "'243.90"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.3131"
$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("D9").Value = "'25.77"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").Value = "'0.07375"
$ws.Range("E10").Value = "  +4.66%  "

$ws.Range("D11").Value = "'0.08075"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "'0.7732"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").Value = "'5.514"
$ws.Range("E13").Value = "  +3.53%  "

$ws.Range("D14").Value = "'1.907.71"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").Value = "'94.34"
$ws.Range("E15").Value = "  +2.10%  "

$ws.Range("D16").Value = "'6.232"
$ws.Range("E16").Value = "  +4.26%  "

$ws.Range("D17").Value = "'29.937.98"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "'14.01"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "'247.53"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("D20").Value = "'0.000007846"
$ws.Range("E20").Value = "  +1.79%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'8.158"
$ws.Range("E21").Value = "  -1.40%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.147.60"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'0.1583"
$ws.Range("E25").Value = "  -4.96%  "

$ws.Range("D26").Value = "'9.474"
$ws.Range("E26").Value = "  +1.51%  "

$ws.Range("D27").Value = "'163.39"
$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").Value = "'18.73"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'2.033"
$ws.Range("E29").Value = "  -1.49%  "

$ws.Range("D30").Value = "'1.432"
$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("D31").Value = "'1.542"
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").Value = "'4.474"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "'0.05576"
$ws.Range("E33").Value = "  -1.87%  "

$ws.Range("D34").Value = "'4.067"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").Value = "'1.244"
$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").Value = "'0.7541"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("D37").Value = "'1.004"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").Value = "'2.684"
$ws.Range("E38").Value = "  +1.71%  "

$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = "  +1.20%  "

$ws.Range("D40").Value = "'2.792"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").Value = "'0.4472"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").Value = "'74.42"
$ws.Range("E42").Value = "  +2.55%  "

$ws.Range("D43").Value = "'1.107.00"
$ws.Range("E43").Value = "  +6.34%  "

$ws.Range("D44").Value = "'5.994"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").Value = "'0.8525"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").Value = "'1.897"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("D48").Value = "'102.66"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.536"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.774"
$ws.Range("E50").Value = "  -1.84%  "

$ws.Range("D51").Value = "'2.990"
$ws.Range("E51").Value = "  +1.97%  "
